# Scheduled-runner style refresh of market/profit data across the Leve
# profit tables (one per crafting job sheet). Only numeric value cells
# in columns H:N are updated; no formulas or formatting are involved.

$wb = $excel.ActiveWorkbook

# Per-sheet, per-row column->value updates. A value of $null means the
# cell should be cleared (removed) rather than written.
$updates = @{
    "ALC" = @{
        92  = @{ H = 642.92;     I = 643.7857;  J = 641.8182;  K = 643.7857;  L = 641.8182;  M = 604.2143;          N = -3137.8182 }
        100 = @{ H = 1177.7391;  I = 926.8125;  J = 1751.2858; K = 926.8125;  L = 1751.2858; M = -385.8125;         N = -2833.2858 }
        137 = @{ H = 41668500;   I = 1692.6923; J = 90911090;  K = 5078.0769; L = 272733270; M = -2528.0769;        N = -272738370 }
    }
    "ARM" = @{
        88 = @{ H = 2366.6667; I = 2660;       J = 2000;      K = 2660;       L = 2000;      M = -2254;             N = -2812 }
        91 = @{ H = 2366.6667; I = 2660;       J = 2000;      K = 2660;       L = 2000;      M = -1256;             N = -4808 }
        97 = @{ H = 554.375;   I = 550.43475;  J = 564.44446; K = 550.43475;  L = 564.44446; M = -54.43475000000001; N = -1556.44446 }
    }
    "BSM" = @{
        80  = @{ H = 423.31818; I = 182.57143; J = 535.6667; K = 182.57143;            L = 535.6667;   M = 815.42857;  N = -2531.6667 }
        83  = @{ H = 423.31818; I = 182.57143; J = 535.6667; K = 912.8571499999999;     L = 2678.3335;  M = 4079.14285; N = -12662.3335 }
        94  = @{ H = 506.3;     I = 255.09091; J = 813.3333; K = 255.09091;             L = 813.3333;   M = 195.90909;  N = -1715.3333 }
        107 = @{ H = 673091.5;  I = 1008259.4; J = 2755.7144; K = 1008259.4;            L = 2755.7144;  M = -1006339.4; N = -6595.7144 }
    }
    "CRP" = @{
        86 = @{ H = 4252.4287; I = 3993.4; J = 4900; K = 3993.4; L = 4900;  M = -2870.4; N = -7146 }
        89 = @{ H = 4252.4287; I = 3993.4; J = 4900; K = 19967;  L = 24500; M = -14351;  N = -35732 }
    }
    "CUL" = @{
        98  = @{ H = 190;        I = 190;       J = 0;         K = 570;       L = 0;        M = 928;                N = $null }
        113 = @{ H = 3036.575;   I = 592.0833;  J = 4084.2144; K = 1776.2499; L = 12252.6432; M = 393.7501;         N = -16592.6432 }
        132 = @{ H = 58824580;   I = 76923860;  J = 1940;      K = 692314740; L = 17460;      M = -692312210;       N = -22520 }
        139 = @{ H = 1910.6923;  I = 1648.7778;  J = 2500;     K = 4946.3334; L = 7500;        M = 193.6665999999996; N = -17780 }
    }
    "GSM" = @{
        70  = @{ H = 4191.758; I = 4072.8333; J = 4334.467; K = 4072.8333; L = 4334.467; M = -3802.8333;  N = -4874.467 }
        73  = @{ H = 4191.758; I = 4072.8333; J = 4334.467; K = 4072.8333; L = 4334.467; M = -3136.8333;  N = -6206.467 }
        107 = @{ H = 562.5417; I = 395.05884; J = 969.2857; K = 395.05884; L = 969.2857; M = 1524.94116;  N = -4809.2857 }
    }
    "LTW" = @{
        93  = @{ H = 1020.86957; I = 670.2857; K = 670.2857; M = 577.7143 }
        133 = @{ H = 59990;      J = 59990;    L = 59990;    N = -65050 }
    }
    "WVR" = @{
        62  = @{ H = 3346;     I = 3152.4614; J = 3975; K = 3152.4614;  L = 3975;   M = -2528.4614; N = -5223 }
        65  = @{ H = 3346;     I = 3152.4614; J = 3975; K = 15762.307;  L = 19875;  M = -12642.307;  N = -26115 }
        136 = @{ H = 10433.519; I = 17575.467; J = 1506.0834; K = 52726.401; L = 4518.2502; M = -50176.401; N = -9618.2502 }
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $cols = $rows[$row]
        foreach ($col in $cols.Keys) {
            $value = $cols[$col]
            $cell = $ws.Range("$col$row")
            if ($null -eq $value) {
                $cell.ClearContents()
            } else {
                $cell.Value = $value
            }
        }
    }
}
